$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.270.11"
$ws.Range("E2").Value = "  +3.78%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.627.87"
$ws.Range("E3").Value = "  +2.13%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "566.96"
$ws.Range("E5").Value = "  +6.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.09"
$ws.Range("E6").Value = "  +3.45%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.37%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.608"
$ws.Range("E8").Value = "  +4.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.650.24"
$ws.Range("E9").Value = "  +2.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.84"
$ws.Range("E10").Value = "  +1.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.105"
$ws.Range("E11").Value = "  +5.93%  "
$ws.Range("E12").Value = "  +6.99%  "
$ws.Range("E13").Value = "  +4.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.097.28"
$ws.Range("E14").Value = "  +2.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.270.67"
$ws.Range("E15").Value = "  +3.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.03"
$ws.Range("E16").Value = "  +6.67%  "
$ws.Range("E17").Value = "  +5.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.652.34"
$ws.Range("E18").Value = "  +3.28%  "
$ws.Range("E19").Value = "  +3.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "342.74"
$ws.Range("E20").Value = "  +2.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.43"
$ws.Range("E21").Value = "  +4.34%  "
$ws.Range("E22").Value = "  +4.18%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.84"
$ws.Range("E24").Value = "  -1.12%  "
$ws.Range("E25").Value = "  +4.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.164"
$ws.Range("E26").Value = "  +3.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.992"
$ws.Range("E27").Value = "  -0.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.39"
$ws.Range("E28").Value = "  +5.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0799"
$ws.Range("E29").Value = "  +10.84%  "
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.70"
$ws.Range("E31").Value = "  +4.73%  "
$ws.Range("E32").Value = "  +5.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "159.79"
$ws.Range("E33").Value = "  +3.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.14"
$ws.Range("E34").Value = "  +1.81%  "
$ws.Range("E35").Value = "  +6.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.892"
$ws.Range("E36").Value = "  +7.92%  "
$ws.Range("E37").Value = "  +5.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.887"
$ws.Range("E38").Value = "  +8.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.59"
$ws.Range("E39").Value = "  +1.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.51"
$ws.Range("E40").Value = "  +7.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "299.00"
$ws.Range("E41").Value = "  +5.96%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.63"
$ws.Range("E42").Value = "  +1.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.995"
$ws.Range("E43").Value = "  -0.48%  "
$ws.Range("E44").Value = "  +4.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.604"
$ws.Range("E45").Value = "  +2.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0544"
$ws.Range("E46").Value = "  +2.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.69"
$ws.Range("E49").Value = "  +0.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0235"
$ws.Range("E50").Value = "  +4.17%  "

# Row 47: Aave -> EnergySwap
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.50"
$ws.Range("E47").Value = "  +6.77%  "

# Row 48: EnergySwap -> Aave
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "128.02"
$ws.Range("E48").Value = "  +16.55%  "

# Row 51: InjectiveProtocol -> RenderToken
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.68"
$ws.Range("E51").Value = "  +7.65%  "
